# Auto-generated: update Price (D) and Volume(1h) (E) columns for cryptos.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.684.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.03%  "
$ws.Range("D3").Value = "'2.718.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -6.33%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'504.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.68%  "
$ws.Range("D6").Value = "'140.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.52%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "'0.529"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.57%  "
$ws.Range("D9").Value = "'2.734.52"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.92%  "
$ws.Range("E10").Value = "  -2.69%  "
$ws.Range("E11").Value = "  +1.97%  "
$ws.Range("E12").Value = "  -3.32%  "
$ws.Range("D13").Value = "'0.126"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("D14").Value = "'3.196.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.20%  "
$ws.Range("D15").Value = "'58.710.67"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.03%  "
$ws.Range("D16").Value = "'21.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.49%  "
$ws.Range("D17").Value = "'2.729.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.22%  "
$ws.Range("D18").Value = "'0.0000135"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.89%  "
$ws.Range("D19").Value = "'4.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.87%  "
$ws.Range("D20").Value = "'10.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.51%  "
$ws.Range("D21").Value = "'340.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.59%  "
$ws.Range("D22").Value = "'6.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.56%  "
$ws.Range("D23").Value = "'0.996"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").Value = "'63.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.28%  "
$ws.Range("D26").Value = "'0.174"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.62%  "
$ws.Range("D27").Value = "'0.426"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.87%  "
$ws.Range("D28").Value = "'0.995"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("D29").Value = "'7.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.30%  "
$ws.Range("D30").Value = "'0.0₃0824"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.31%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("D32").Value = "'19.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.52%  "
$ws.Range("E33").Value = "  -4.70%  "
$ws.Range("D34").Value = "'150.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.92%  "
$ws.Range("D35").Value = "'4.20"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.94%  "
$ws.Range("D36").Value = "'5.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.96%  "
$ws.Range("D37").Value = "'0.945"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.79%  "
$ws.Range("E38").Value = "  -6.37%  "
$ws.Range("D39").Value = "'36.26"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.33%  "
$ws.Range("D40").Value = "'3.57"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.63%  "
$ws.Range("E41").Value = "  -7.39%  "
$ws.Range("D42").Value = "'2.193.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.92%  "
$ws.Range("D43").Value = "'0.0560"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.90%  "
$ws.Range("D44").Value = "'0.997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").Value = "'0.600"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.93%  "
$ws.Range("D46").Value = "'18.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.31%  "
$ws.Range("D47").Value = "'4.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.33%  "
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("E49").Value = "  -3.33%  "
$ws.Range("D50").Value = "'0.0888"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.95%  "
$ws.Range("D51").Value = "'18.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.31%  "
